$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.799257602882035, 1.5021707826232191, -1.1517886662636974, 2.656790039974177)"
$ws.Range("C2").Value = "NIG(1.3972564691830003, 1.1160143088384082, 4.952284557658946, 5.849337751043771)"
$ws.Range("D2").Value = "NIG(0.8378493481521685, 0.5587806432168014, 1.3801038347749603, 3.3503922326746287)"
$ws.Range("E2").Value = "NCT(2.4556167987130504, 1.51701196430324, -0.009623588115898571, 4.483896015283383)"
